$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A36").Value = 2943549
$ws.Range("B36").Value = "Pril Power blue 650 ml"
$ws.Range("C36").Value = 12
$ws.Range("D36").Value = 168
